$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Cells.Item(1,1).Value = "ISIN"
$ws.Cells.Item(1,2).Value = "Stock Name"
$ws.Cells.Item(1,3).Value = "Mutual Fund"
$ws.Cells.Item(1,4).Value = "Status"
$ws.Cells.Item(1,5).Value = "Jan_2026"
$ws.Cells.Item(1,6).Value = "Dec_2025"
$ws.Cells.Item(1,7).Value = "Oct_2025"
$ws.Cells.Item(1,8).Value = "MoM"
$ws.Cells.Item(1,9).Value = "QoQ"

# Apply the existing header style (bold, border, centered) to the new I1 header cell
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows ---
# Row 2: INE423A01024 - Adani Enterprises Limited
$ws.Cells.Item(2,1).Value = "INE423A01024"
$ws.Cells.Item(2,2).Value = "Adani Enterprises Limited"
$ws.Cells.Item(2,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(2,4).Value = "Reducing"
$ws.Cells.Item(2,5).Value = 8.711195999999999
$ws.Cells.Item(2,6).Value = 9.043240000000001
$ws.Cells.Item(2,7).Value = 4.504588
$ws.Cells.Item(2,8).Value = -0.3320440000000016
$ws.Cells.Item(2,9).Value = 4.206607999999999

# Row 3: INE364U01010 - Adani Green Energy Limited
$ws.Cells.Item(3,1).Value = "INE364U01010"
$ws.Cells.Item(3,2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(3,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(3,4).Value = "Reducing"
$ws.Cells.Item(3,5).Value = 7.992836
$ws.Cells.Item(3,6).Value = 8.914787
$ws.Cells.Item(3,7).Value = 4.663532
$ws.Cells.Item(3,8).Value = -0.9219510000000009
$ws.Cells.Item(3,9).Value = 3.329304

# Row 4: INE202B01038 - Piramal Finance Ltd
$ws.Cells.Item(4,1).Value = "INE202B01038"
$ws.Cells.Item(4,2).Value = "Piramal Finance Ltd"
$ws.Cells.Item(4,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(4,4).Value = "Adding"
$ws.Cells.Item(4,5).Value = 7.951795
$ws.Cells.Item(4,6).Value = 7.017373
$ws.Cells.Item(4,7).Value = 8.325765000000001
$ws.Cells.Item(4,8).Value = 0.9344219999999996
$ws.Cells.Item(4,9).Value = -0.3739700000000008

# Row 5: INE917I01010 - Bajaj Auto Limited
$ws.Cells.Item(5,1).Value = "INE917I01010"
$ws.Cells.Item(5,2).Value = "Bajaj Auto Limited"
$ws.Cells.Item(5,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(5,4).Value = "Adding Consistently"
$ws.Cells.Item(5,5).Value = 7.515284
$ws.Cells.Item(5,6).Value = 6.851212
$ws.Cells.Item(5,7).Value = 6.399446
$ws.Cells.Item(5,8).Value = 0.664072
$ws.Cells.Item(5,9).Value = 1.115838

# Row 6: INE931S01010 - Adani Energy Solutions Limited
$ws.Cells.Item(6,1).Value = "INE931S01010"
$ws.Cells.Item(6,2).Value = "Adani Energy Solutions Limited"
$ws.Cells.Item(6,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(6,4).Value = "Reducing"
$ws.Cells.Item(6,5).Value = 6.68311
$ws.Cells.Item(6,6).Value = 7.18563
$ws.Cells.Item(6,7).Value = 4.032672
$ws.Cells.Item(6,8).Value = -0.5025199999999996
$ws.Cells.Item(6,9).Value = 2.650438

# Row 7: INE726G01019 - ICICI Prudential Life Insurance Co Ltd
$ws.Cells.Item(7,1).Value = "INE726G01019"
$ws.Cells.Item(7,2).Value = "ICICI Prudential Life Insurance Co Ltd"
$ws.Cells.Item(7,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(7,4).Value = "Fresh Entry"
$ws.Cells.Item(7,5).Value = 5.858935
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 0
$ws.Cells.Item(7,8).Value = 5.858935
$ws.Cells.Item(7,9).Value = 5.858935

# Row 8: INE090A01021 - ICICI Bank Limited
$ws.Cells.Item(8,1).Value = "INE090A01021"
$ws.Cells.Item(8,2).Value = "ICICI Bank Limited"
$ws.Cells.Item(8,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(8,4).Value = "Fresh Entry"
$ws.Cells.Item(8,5).Value = 5.626443
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = 0
$ws.Cells.Item(8,8).Value = 5.626443
$ws.Cells.Item(8,9).Value = 5.626443

# Row 9: INE406A01037 - Aurobindo Pharma Limited
$ws.Cells.Item(9,1).Value = "INE406A01037"
$ws.Cells.Item(9,2).Value = "Aurobindo Pharma Limited"
$ws.Cells.Item(9,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(9,4).Value = "Adding Consistently"
$ws.Cells.Item(9,5).Value = 5.354146
$ws.Cells.Item(9,6).Value = 4.911449
$ws.Cells.Item(9,7).Value = 4.64032
$ws.Cells.Item(9,8).Value = 0.4426969999999999
$ws.Cells.Item(9,9).Value = 0.7138260000000001

# Row 10: INE016A01026 - Dabur India Limited
$ws.Cells.Item(10,1).Value = "INE016A01026"
$ws.Cells.Item(10,2).Value = "Dabur India Limited"
$ws.Cells.Item(10,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(10,4).Value = "Fresh Entry"
$ws.Cells.Item(10,5).Value = 4.123527
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = 4.123527
$ws.Cells.Item(10,9).Value = 4.123527

# Row 11: INE237A01036 - Kotak Mahindra Bank Limited
$ws.Cells.Item(11,1).Value = "INE237A01036"
$ws.Cells.Item(11,2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(11,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(11,4).Value = "Fresh Entry"
$ws.Cells.Item(11,5).Value = 4.089286
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = 0
$ws.Cells.Item(11,8).Value = 4.089286
$ws.Cells.Item(11,9).Value = 4.089286

# Row 12: INE180C01042 - Capri Global Capital Limited
$ws.Cells.Item(12,1).Value = "INE180C01042"
$ws.Cells.Item(12,2).Value = "Capri Global Capital Limited"
$ws.Cells.Item(12,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(12,4).Value = "Adding"
$ws.Cells.Item(12,5).Value = 2.77975
$ws.Cells.Item(12,6).Value = 2.708321
$ws.Cells.Item(12,7).Value = 2.95139
$ws.Cells.Item(12,8).Value = 0.07142899999999974
$ws.Cells.Item(12,9).Value = -0.17164

# Row 13: INE040A01034 - HDFC Bank Limited
$ws.Cells.Item(13,1).Value = "INE040A01034"
$ws.Cells.Item(13,2).Value = "HDFC Bank Limited"
$ws.Cells.Item(13,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(13,4).Value = "Fresh Entry"
$ws.Cells.Item(13,5).Value = 2.50462
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 0
$ws.Cells.Item(13,8).Value = 2.50462
$ws.Cells.Item(13,9).Value = 2.50462

# Row 14: INE860A01027 - HCL Technologies Limited
$ws.Cells.Item(14,1).Value = "INE860A01027"
$ws.Cells.Item(14,2).Value = "HCL Technologies Limited"
$ws.Cells.Item(14,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(14,4).Value = "Complete Exit"
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = 1.902312
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = -1.902312
$ws.Cells.Item(14,9).Value = 0

# Row 15: INE814H01029 - Adani Power Limited
$ws.Cells.Item(15,1).Value = "INE814H01029"
$ws.Cells.Item(15,2).Value = "Adani Power Limited"
$ws.Cells.Item(15,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(15,4).Value = "Complete Exit"
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = 5.695908
$ws.Cells.Item(15,8).Value = 0
$ws.Cells.Item(15,9).Value = -5.695908

# Row 16: INE669C01036 - Tech Mahindra Limited
$ws.Cells.Item(16,1).Value = "INE669C01036"
$ws.Cells.Item(16,2).Value = "Tech Mahindra Limited"
$ws.Cells.Item(16,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(16,4).Value = "Complete Exit"
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 3.371679
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(16,8).Value = -3.371679
$ws.Cells.Item(16,9).Value = 0

# Row 17: INE775A01035 - Samvardhana Motherson International Ltd
$ws.Cells.Item(17,1).Value = "INE775A01035"
$ws.Cells.Item(17,2).Value = "Samvardhana Motherson International Ltd"
$ws.Cells.Item(17,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(17,4).Value = "Complete Exit"
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = 8.769418999999999
$ws.Cells.Item(17,7).Value = 0
$ws.Cells.Item(17,8).Value = -8.769418999999999
$ws.Cells.Item(17,9).Value = 0

# Row 18: INE768C01028 - Zydus Wellness Ltd
$ws.Cells.Item(18,1).Value = "INE768C01028"
$ws.Cells.Item(18,2).Value = "Zydus Wellness Ltd"
$ws.Cells.Item(18,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(18,4).Value = "Complete Exit"
$ws.Cells.Item(18,5).Value = 0
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = 7.530469
$ws.Cells.Item(18,8).Value = 0
$ws.Cells.Item(18,9).Value = -7.530469

# Row 19: INE009A01021 - Infosys Limited
$ws.Cells.Item(19,1).Value = "INE009A01021"
$ws.Cells.Item(19,2).Value = "Infosys Limited"
$ws.Cells.Item(19,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(19,4).Value = "Complete Exit"
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = 0
$ws.Cells.Item(19,7).Value = 3.959041
$ws.Cells.Item(19,8).Value = 0
$ws.Cells.Item(19,9).Value = -3.959041

# Row 20: INE484J01027 - Godrej Properties Limited
$ws.Cells.Item(20,1).Value = "INE484J01027"
$ws.Cells.Item(20,2).Value = "Godrej Properties Limited"
$ws.Cells.Item(20,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(20,4).Value = "Complete Exit"
$ws.Cells.Item(20,5).Value = 0
$ws.Cells.Item(20,6).Value = 3.130396
$ws.Cells.Item(20,7).Value = 0.63528
$ws.Cells.Item(20,8).Value = -3.130396
$ws.Cells.Item(20,9).Value = -0.63528

# Row 21: INE00H001014 - SWIGGY LIMITED
$ws.Cells.Item(21,1).Value = "INE00H001014"
$ws.Cells.Item(21,2).Value = "SWIGGY LIMITED"
$ws.Cells.Item(21,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(21,4).Value = "Complete Exit"
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 2.659029
$ws.Cells.Item(21,8).Value = 0
$ws.Cells.Item(21,9).Value = -2.659029

# Row 22: INE245A01021 - Tata Power Company Limited
$ws.Cells.Item(22,1).Value = "INE245A01021"
$ws.Cells.Item(22,2).Value = "Tata Power Company Limited"
$ws.Cells.Item(22,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(22,4).Value = "Complete Exit"
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = 6.230446
$ws.Cells.Item(22,8).Value = 0
$ws.Cells.Item(22,9).Value = -6.230446

# Row 23: INE237A01028 - Kotak Mahindra Bank Limited
$ws.Cells.Item(23,1).Value = "INE237A01028"
$ws.Cells.Item(23,2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(23,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(23,4).Value = "Complete Exit"
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 4.131909
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(23,8).Value = -4.131909
$ws.Cells.Item(23,9).Value = 0

# Row 24: INE192A01025 - Tata Consumer Products Ltd
$ws.Cells.Item(24,1).Value = "INE192A01025"
$ws.Cells.Item(24,2).Value = "Tata Consumer Products Ltd"
$ws.Cells.Item(24,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(24,4).Value = "Complete Exit"
$ws.Cells.Item(24,5).Value = 0
$ws.Cells.Item(24,6).Value = 0
$ws.Cells.Item(24,7).Value = 3.653407
$ws.Cells.Item(24,8).Value = 0
$ws.Cells.Item(24,9).Value = -3.653407

# Row 25: INE047A01021 - Grasim Industries Ltd
$ws.Cells.Item(25,1).Value = "INE047A01021"
$ws.Cells.Item(25,2).Value = "Grasim Industries Ltd"
$ws.Cells.Item(25,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(25,4).Value = "Complete Exit"
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(25,6).Value = 0
$ws.Cells.Item(25,7).Value = 2.985991
$ws.Cells.Item(25,8).Value = 0
$ws.Cells.Item(25,9).Value = -2.985991

# Row 26: INE019A01038 - JSW Steel Limited
$ws.Cells.Item(26,1).Value = "INE019A01038"
$ws.Cells.Item(26,2).Value = "JSW Steel Limited"
$ws.Cells.Item(26,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(26,4).Value = "Complete Exit"
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(26,6).Value = 0
$ws.Cells.Item(26,7).Value = 1.528663
$ws.Cells.Item(26,8).Value = 0
$ws.Cells.Item(26,9).Value = -1.528663

# Row 27: INE271C01023 - DLF Limited
$ws.Cells.Item(27,1).Value = "INE271C01023"
$ws.Cells.Item(27,2).Value = "DLF Limited"
$ws.Cells.Item(27,3).Value = "quant ESG Integration Strategy Fund"
$ws.Cells.Item(27,4).Value = "Complete Exit"
$ws.Cells.Item(27,5).Value = 0
$ws.Cells.Item(27,6).Value = 7.00138
$ws.Cells.Item(27,7).Value = 9.556981
$ws.Cells.Item(27,8).Value = -7.00138
$ws.Cells.Item(27,9).Value = -9.556981
